$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C, rows 2-23,
# from serial 45175 (2023-09-06) to serial 45183 (2023-09-14).
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
